$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new M:N table
$ws.Range("M1").Value = "V1"
$ws.Range("N1").Value = "V2"

# Data rows (row 8 is intentionally skipped - no data there)
$data = @{
    2  = @(0.05, 0.06)
    3  = @(0.1, 0.14000000000000001)
    4  = @(0.15, 0.2)
    5  = @(0.2, 0.28000000000000003)
    6  = @(0.25, 0.36)
    7  = @(0.3, 0.4)
    9  = @(0.35, 0.5)
    10 = @(0.4, 0.6)
    11 = @(0.45, 0.7)
    12 = @(0.5, 0.75)
    13 = @(0.75, 1.4)
    14 = @(1, 2)
    15 = @(2, 5)
    16 = @(3, 9)
    17 = @(4, 17.5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 13).Value = $vals[0]
    $ws.Cells.Item($row, 14).Value = $vals[1]
}

# Update selection to match the new active range
$ws.Range("M2:N18").Select()
